$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 2.31244962562356
$ws.Range("D8").Value = 0.0948663428343971
$ws.Range("E8").Value = 0.577820197044335
$ws.Range("F8").Value = 0.756521739130435
$ws.Range("G8").Value = 0.978107689448793
